$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamps = @{
    2  = "2025-11-03T00:10:28.840636"
    3  = "2025-11-03T00:10:28.840636"
    4  = "2025-11-03T00:10:28.840636"
    5  = "2025-11-03T00:10:28.840636"
    6  = "2025-11-03T00:10:28.840636"
    7  = "2025-11-03T00:10:28.840636"
    8  = "2025-11-03T00:10:28.841630"
    9  = "2025-11-03T00:10:28.842639"
    10 = "2025-11-03T00:10:28.842639"
    11 = "2025-11-03T00:10:28.843636"
    12 = "2025-11-03T00:10:28.843636"
    13 = "2025-11-03T00:10:28.844636"
    14 = "2025-11-03T00:10:28.844636"
    15 = "2025-11-03T00:10:28.844636"
    16 = "2025-11-03T00:10:28.845638"
    17 = "2025-11-03T00:10:28.845638"
    18 = "2025-11-03T00:10:28.845638"
    19 = "2025-11-03T00:10:28.845638"
    20 = "2025-11-03T00:10:28.845638"
    21 = "2025-11-03T00:10:28.846635"
    22 = "2025-11-03T00:10:28.846635"
    23 = "2025-11-03T00:10:28.846635"
    24 = "2025-11-03T00:10:28.846635"
    25 = "2025-11-03T00:10:28.846635"
    26 = "2025-11-03T00:10:28.847629"
    27 = "2025-11-03T00:10:28.847629"
    28 = "2025-11-03T00:10:28.847629"
    29 = "2025-11-03T00:10:28.847629"
}

foreach ($row in $timestamps.Keys) {
    $ws.Range("Z$row").Value = $timestamps[$row]
}
